$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G273").Value = 1291
$ws.Range("F289").Value = 62799
$ws.Range("F299").Value = 64522
$ws.Range("G299").Value = 6799
$ws.Range("F300").Value = 71394
$ws.Range("G300").Value = 7030
$ws.Range("F306").Value = 71013
$ws.Range("G306").Value = 7159
$ws.Range("F307").Value = 75142
$ws.Range("G307").Value = 6412
$ws.Range("F308").Value = 15623
$ws.Range("G308").Value = 1084
$ws.Range("F309").Value = 74466
$ws.Range("G309").Value = 5279
$ws.Range("F310").Value = 75184
$ws.Range("G310").Value = 3914
$ws.Range("F311").Value = 62148
$ws.Range("G311").Value = 1971
$ws.Range("F313").Value = 71508
$ws.Range("G313").Value = 3212
$ws.Range("F314").Value = 63413
$ws.Range("F315").Value = 55842
$ws.Range("F316").Value = 49219
$ws.Range("G316").Value = 2227
$ws.Range("F317").Value = 61694
$ws.Range("G317").Value = 2114
$ws.Range("F318").Value = 49845
$ws.Range("G318").Value = 1198
$ws.Range("F320").Value = 76243
$ws.Range("G320").Value = 3653
$ws.Range("F321").Value = 90621
$ws.Range("G321").Value = 2795
$ws.Range("F322").Value = 106338
$ws.Range("F323").Value = 213065
$ws.Range("G323").Value = 3162
$ws.Range("F324").Value = 232967
$ws.Range("G324").Value = 2664
$ws.Range("F325").Value = 754808
$ws.Range("G325").Value = 6371
$ws.Range("F329").Value = 89151
$ws.Range("G329").Value = 1810
$ws.Range("F331").Value = 150383
$ws.Range("G331").Value = 2595
$ws.Range("F332").Value = 427724
$ws.Range("G332").Value = 4182
$ws.Range("F333").Value = 259813
$ws.Range("G333").Value = 2801
$ws.Range("F334").Value = 203633
$ws.Range("G334").Value = 3389
$ws.Range("F335").Value = 129871
$ws.Range("G335").Value = 2927
$ws.Range("F336").Value = 100587
$ws.Range("G336").Value = 3192
$ws.Range("F337").Value = 102289
$ws.Range("G337").Value = 2896
$ws.Range("F338").Value = 216674
$ws.Range("G338").Value = 3050
$ws.Range("F339").Value = 632067
$ws.Range("G339").Value = 5392
$ws.Range("F340").Value = 376608
$ws.Range("G340").Value = 3211
$ws.Range("F341").Value = 297564
$ws.Range("G341").Value = 3684
$ws.Range("F342").Value = 186534
$ws.Range("G342").Value = 3160
$ws.Range("F343").Value = 124639
$ws.Range("G343").Value = 2777
$ws.Range("F344").Value = 128423
$ws.Range("G344").Value = 2385
$ws.Range("F345").Value = 271425
$ws.Range("G345").Value = 3204
$ws.Range("F346").Value = 624815
$ws.Range("G346").Value = 4388
$ws.Range("F347").Value = 320592
$ws.Range("G347").Value = 2722
$ws.Range("F348").Value = 218334
$ws.Range("G348").Value = 2968
